$d = $word.ActiveDocument

$pairs = @(
    @("2024-03-01 Friday", "2024-03-02 Saturday"),
    @("66+26=92", "5+78=83"),
    @("50-23=27", "44+27=71"),
    @("62-54=8", "81-48=33"),
    @("60+2=62", "39-39=0"),
    @("47+2=49", "22-6=16"),
    @("15+84=99", "25+50=75"),
    @("28-2=26", "23+75=98"),
    @("82+14=96", "38+14=52"),
    @("94-21=73", "53-34=19"),
    @("98-7=91", "31+31=62"),
    @("70+1=71", "83-43=40"),
    @("55+9=64", "98-51=47"),
    @("58+8=66", "55-23=32"),
    @("66+22=88", "5+78=83"),
    @("60-9=51", "10+28=38"),
    @("99-3=96", "44+46=90"),
    @("34+38=72", "35-3=32"),
    @("19+32=51", "99-42=57"),
    @("84-32=52", "60-6=54"),
    @("11-1=10", "39+20=59"),
    @("44+4=48", "17+38=55"),
    @("22+75=97", "18+69=87"),
    @("10-8=2", "11+42=53"),
    @("25+35=60", "26-21=5"),
    @("88-27=61", "67-44=23"),
    @("57+24=81", "25+0=25"),
    @("30-8=22", "71+4=75"),
    @("69+14=83", "64-53=11"),
    @("9+72=81", "55-44=11"),
    @("8+37=45", "19-7=12"),
    @("85-66=19", "66-30=36"),
    @("27-10=17", "8+77=85"),
    @("57-11=46", "19+47=66"),
    @("27+0=27", "15+12=27"),
    @("57+25=82", "69+27=96"),
    @("26+53=79", "6-2=4"),
    @("20+25=45", "33-11=22"),
    @("69-62=7", "20+70=90"),
    @("48+17=65", "6+82=88"),
    @("94-33=61", "60-42=18"),
    @("70-52=18", "38+58=96"),
    @("37+0=37", "40+9=49"),
    @("82-57=25", "36+1=37"),
    @("19+25=44", "10+82=92"),
    @("89+3=92", "32+59=91"),
    @("9+51=60", "73-17=56"),
    @("47+4=51", "48+15=63"),
    @("71-20=51", "15+69=84"),
    @("83-13=70", "16+1=17"),
    @("67-57=10", "36+11=47"),
    @("53+13=66", "14+23=37"),
    @("25+69=94", "0+13=13"),
    @("67+13=80", "23+52=75"),
    @("24+9=33", "11+11=22"),
    @("77-59=18", "0+27=27"),
    @("50+38=88", "38+7=45"),
    @("59+27=86", "59+12=71"),
    @("70+7=77", "75-48=27"),
    @("40+2=42", "9+70=79"),
    @("93-24=69", "68-24=44"),
    @("2+38=40", "62+6=68"),
    @("82-79=3", "87+0=87"),
    @("5+43=48", "61-20=41"),
    @("71+7=78", "25-18=7"),
    @("5+88=93", "60-0=60"),
    @("81+18=99", "28-20=8"),
    @("20+66=86", "32+38=70"),
    @("93-52=41", "94-57=37"),
    @("96-88=8", "65+11=76"),
    @("15+37=52", "29+14=43"),
    @("64+19=83", "49-39=10"),
    @("43+50=93", "51+14=65"),
    @("15+43=58", "12+35=47"),
    @("59+1=60", "3+67=70"),
    @("93-75=18", "78-22=56"),
    @("48-12=36", "42+22=64"),
    @("89-28=61", "20+73=93"),
    @("68+0=68", "81-75=6"),
    @("33+17=50", "15+13=28"),
    @("90-10=80", "61+8=69"),
    @("2-2=0", "80-49=31"),
    @("71-59=12", "2+53=55"),
    @("38-7=31", "10+79=89"),
    @("32+32=64", "94-29=65"),
    @("46-27=19", "49+22=71"),
    @("68-55=13", "69+5=74"),
    @("20+51=71", "51+31=82"),
    @("86-74=12", "36+17=53"),
    @("90-67=23", "93-15=78"),
    @("19+18=37", "50-42=8"),
    @("73-33=40", "43-15=28"),
    @("69-29=40", "38+13=51"),
    @("61-60=1", "60-35=25"),
    @("32+54=86", "18+54=72"),
    @("52+31=83", "35-16=19"),
    @("54-21=33", "84-58=26"),
    @("68-29=39", "64-47=17"),
    @("10+56=66", "46+0=46"),
    @("63+9=72", "21+41=62"),
    @("51+40=91", "71-44=27")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
